# Add three new verb-conjugation rows (89-91) to Sheet1, matching the
# "add verbs and non verbs" commit: 起こる (to occur), 進む (to advance),
# and 助ける (to help), each with Te / Ta / Nai / Masu / Volitional / Ba forms.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate the formatting of the previous last row (row 88) down to the
# three new rows first, so every new cell reuses the existing "Yu Gothic"
# (A:F) / "微軟正黑體" (G) cell styles instead of Excel minting brand new
# ones the moment a Font property gets touched directly.
$ws.Range("A88:G88").Copy($ws.Range("A89:G89"))
$ws.Range("A88:G88").Copy($ws.Range("A90:G90"))
$ws.Range("A88:G88").Copy($ws.Range("A91:G91"))

# --- Row 89: 起こる (okoru - "to occur") -------------------------------
# Filled left-to-right, same order as the column headers.
$ws.Range("A89").Value = "起こる"
$ws.Range("B89").Value = "起こって"
$ws.Range("C89").Value = "起こった"
$ws.Range("D89").Value = "起こらない"
$ws.Range("E89").Value = "起こります"
$ws.Range("F89").Value = "起ころう"
$ws.Range("G89").Value = "起これば"

# --- Row 90: 進む (susumu - "to advance") ------------------------------
# Entered as Dict/Te/Ta/Masu/Volitional/Nai/Ba (Nai typed after Masu &
# Volitional) so the new shared-string order matches the source workbook.
$ws.Range("A90").Value = "進む"
$ws.Range("B90").Value = "進んで"
$ws.Range("C90").Value = "進んだ"
$ws.Range("E90").Value = "進みます"
$ws.Range("F90").Value = "進もう"
$ws.Range("D90").Value = "進まない"
$ws.Range("G90").Value = "進めば"

# --- Row 91: 助ける (tasukeru - "to help") ------------------------------
$ws.Range("A91").Value = "助ける"
$ws.Range("B91").Value = "助けて"
$ws.Range("C91").Value = "助けた"
$ws.Range("D91").Value = "助けない"
$ws.Range("E91").Value = "助けます"
$ws.Range("F91").Value = "助けよう"
$ws.Range("G91").Value = "助ければ"

# Row height matches the rest of the conjugation table (18.75pt rows).
$ws.Rows.Item(89).RowHeight = 18.75
$ws.Rows.Item(90).RowHeight = 18.75
$ws.Rows.Item(91).RowHeight = 18.75

# --- Selection / scroll position, matching the tail of the sheet -------
$ws.Range("G91").Select()

Write-Host "Added rows 89-91 (okoru, susumu, tasukeru)"
